$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark from the first paragraph (Q1).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Merge the two runs "Please provide a " + "name and a way to contact you"
#    into a single run with the combined text (Find/Replace collapses them).
$q4Para = $d.Paragraphs.Item(14)
[void]$q4Para.Range.Find.Execute("Please provide a name and a way to contact you", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Please provide a name and a way to contact you", 2)

# 3. Delete the four "Name / Company / City/Town / Email Address" paragraphs.
$startPara = $d.Paragraphs.Item(15)
$endPara = $d.Paragraphs.Item(18)
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# 4. Re-add the _GoBack bookmark (collapsed) at the start of the final
#    paragraph, which now holds the drawing.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmRange = $lastPara.Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)
